$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell value updates from the latest crypto price scrape.
# D-column price cells that would otherwise be auto-parsed as
# numbers are forced to remain plain text (matching the source data,
# which stores prices/volumes as text strings), by setting the
# cell number format to Text ("@") before assigning the value.

$ws.Range("D2").Value = "36.570.93"
$ws.Range("E2").Value = "  -1.40%  "
$ws.Range("D3").Value = "2.063.39"
$ws.Range("E3").Value = "  +0.95%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.00"
$ws.Range("E5").Value = "  -2.64%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.675"
$ws.Range("E6").Value = "  +1.86%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "52.34"
$ws.Range("E8").Value = "  -6.89%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "58.83"
$ws.Range("E9").Value = "  -1.82%  "
$ws.Range("E10").Value = "  -6.13%  "
$ws.Range("E11").Value = "  -3.66%  "
$ws.Range("E12").Value = "  -1.34%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.885"
$ws.Range("E13").Value = "  -1.53%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.44"
$ws.Range("E14").Value = "  -9.73%  "
$ws.Range("D15").Value = "2.366.45"
$ws.Range("E15").Value = "  +1.08%  "
$ws.Range("E16").Value = "  -5.26%  "
$ws.Range("D17").Value = "2.072.54"
$ws.Range("E17").Value = "  +1.36%  "
$ws.Range("D18").Value = "36.474.41"
$ws.Range("E18").Value = "  -1.69%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.28"
$ws.Range("E19").Value = "  -13.48%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.62"
$ws.Range("E20").Value = "  -4.41%  "
$ws.Range("D21").Value = "0.0₃0859"
$ws.Range("E21").Value = "  -3.56%  "
$ws.Range("E22").Value = "  -2.77%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.32"
$ws.Range("E23").Value = "  -0.70%  "
$ws.Range("E24").Value = "  -0.08%  "
$ws.Range("E25").Value = "  -5.09%  "
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.23"
$ws.Range("E26").Value = "  -3.32%  "
$ws.Range("B27").Value = "PancakeSwap"
$ws.Range("C27").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.12"
$ws.Range("E27").Value = "  -2.42%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "163.02"
$ws.Range("E28").Value = "  -4.65%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "20.25"
$ws.Range("E29").Value = "  +0.69%  "
$ws.Range("E30").Value = "  -0.98%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.05"
$ws.Range("E31").Value = "  -1.34%  "
$ws.Range("E32").Value = "  -3.97%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.55"
$ws.Range("E33").Value = "  -2.25%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0594"
$ws.Range("E34").Value = "  -4.44%  "
$ws.Range("E35").Value = "  +0.04%  "
$ws.Range("B36").Value = "WEMIXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.84"
$ws.Range("E36").Value = "  -1.50%  "
$ws.Range("B37").Value = "LidoDAOToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.26"
$ws.Range("E37").Value = "  -0.45%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0810"
$ws.Range("E38").Value = "  -7.48%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.24"
$ws.Range("E39").Value = "  -6.48%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.84"
$ws.Range("E40").Value = "  -6.10%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.88"
$ws.Range("E41").Value = "  -6.16%  "
$ws.Range("E42").Value = "  -3.82%  "
$ws.Range("E43").Value = "  -2.73%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0928"
$ws.Range("E44").Value = "  -6.43%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "93.61"
$ws.Range("E45").Value = "  -5.32%  "
$ws.Range("D46").Value = "1.384.63"
$ws.Range("E46").Value = "  +7.74%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "15.53"
$ws.Range("E47").Value = "  -9.80%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.25"
$ws.Range("E48").Value = "  +6.49%  "
$ws.Range("E49").Value = "  -4.07%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.84"
$ws.Range("E50").Value = "  -0.45%  "
$ws.Range("D51").Value = "2.250.69"
$ws.Range("E51").Value = "  +1.14%  "
